# "Updated cryptos list on Sat Jul 22 11:44:56 UTC 2023 with GitHub Actions"
#
# Refreshes the Price (D) / Volume(1h) (E) columns of the cryptos table on
# Sheet1 with a newly scraped snapshot, rows 2-51. A handful of coins that
# swapped market-cap rank also get their Coin (B) / Link (C) cells updated
# so the row order matches the new ranking (WrappedEther<->Polygon,
# PancakeSwap<->Filecoin, Aave<->FraxShare).
#
# Several Price cells are plain-text decimal numbers (e.g. "1.001",
# "13.80") whose exact digits -- including trailing zeros -- must be kept
# verbatim, just like the original sheet stores them as text. A bare
# Range.Value assignment lets Excel's type inference coerce a string such
# as "1.000" into the number 1, so for any new value that looks numeric we
# momentarily force the cell to Text format, assign it, then clear the
# formatting again (ClearFormats leaves the stored text value untouched
# while restoring the cell's original/default style).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '29.893.58'
$ws.Range('E2').Value = '  +0.40%  '
$ws.Range('D3').Value = '1.889.09'
$ws.Range('E3').Value = '  +0.08%  '
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '1.001'
$ws.Range('D4').ClearFormats()
$ws.Range('E4').Value = '  -0.01%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '0.7743'
$ws.Range('D5').ClearFormats()
$ws.Range('E5').Value = '  +0.18%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '242.89'
$ws.Range('D6').ClearFormats()
$ws.Range('E6').Value = '  -0.57%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '1.000'
$ws.Range('D7').ClearFormats()
$ws.Range('E7').Value = '  -0.02%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.3118'
$ws.Range('D8').ClearFormats()
$ws.Range('E8').Value = '  -0.34%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '25.65'
$ws.Range('D9').ClearFormats()
$ws.Range('E9').Value = '  +1.67%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.07188'
$ws.Range('D10').ClearFormats()
$ws.Range('E10').Value = '  -0.45%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.08615'
$ws.Range('D11').ClearFormats()
$ws.Range('E11').Value = '  +6.55%  '
$ws.Range('B12').Value = 'Polygon'
$ws.Range('C12').Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.7644'
$ws.Range('D12').ClearFormats()
$ws.Range('E12').Value = '  -0.03%  '
$ws.Range('B13').Value = 'WrappedEther'
$ws.Range('C13').Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range('D13').Value = '1.934.76'
$ws.Range('E13').Value = '  +0.71%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '5.375'
$ws.Range('D14').ClearFormats()
$ws.Range('E14').Value = '  -1.81%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '93.85'
$ws.Range('D15').ClearFormats()
$ws.Range('E15').Value = '  +1.79%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '6.193'
$ws.Range('D16').ClearFormats()
$ws.Range('E16').Value = '  +0.36%  '
$ws.Range('D17').Value = '29.949.67'
$ws.Range('E17').Value = '  +0.59%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '13.80'
$ws.Range('D18').ClearFormats()
$ws.Range('E18').Value = '  -0.61%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '244.63'
$ws.Range('D19').ClearFormats()
$ws.Range('E19').Value = '  +0.86%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '0.000007834'
$ws.Range('D20').ClearFormats()
$ws.Range('E20').Value = '  +0.81%  '
$ws.Range('D21').Value = '2.194.00'
$ws.Range('E21').Value = '  +2.92%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '0.9981'
$ws.Range('D22').ClearFormats()
$ws.Range('E22').Value = '  -0.30%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '8.023'
$ws.Range('D23').ClearFormats()
$ws.Range('E23').Value = '  -1.66%  '
$ws.Range('E24').Value = '  -0.01%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '0.1641'
$ws.Range('D25').ClearFormats()
$ws.Range('E25').Value = '  +3.64%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '9.385'
$ws.Range('D26').ClearFormats()
$ws.Range('E26').Value = '  -0.38%  '
$ws.Range('E27').Value = '  -0.10%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '18.77'
$ws.Range('D28').ClearFormats()
$ws.Range('E28').Value = '  +0.23%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '2.039'
$ws.Range('D29').ClearFormats()
$ws.Range('E29').Value = '  +0.16%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '1.443'
$ws.Range('D30').ClearFormats()
$ws.Range('E30').Value = '  +0.23%  '
$ws.Range('B31').Value = 'Filecoin'
$ws.Range('C31').Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '4.535'
$ws.Range('D31').ClearFormats()
$ws.Range('E31').Value = '  +1.80%  '
$ws.Range('B32').Value = 'PancakeSwap'
$ws.Range('C32').Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '1.533'
$ws.Range('D32').ClearFormats()
$ws.Range('E32').Value = '  -0.80%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '4.107'
$ws.Range('D33').ClearFormats()
$ws.Range('E33').Value = '  +0.72%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '0.05432'
$ws.Range('D34').ClearFormats()
$ws.Range('E34').Value = '  -1.14%  '
$ws.Range('E35').Value = '  -1.10%  '
$ws.Range('E36').Value = '  -0.70%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '1.004'
$ws.Range('D37').ClearFormats()
$ws.Range('E37').Value = '  +0.54%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '2.697'
$ws.Range('D38').ClearFormats()
$ws.Range('E38').Value = '  +2.34%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.01966'
$ws.Range('D39').ClearFormats()
$ws.Range('E39').Value = '  +2.70%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '2.783'
$ws.Range('D40').ClearFormats()
$ws.Range('E40').Value = '  +0.05%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.4470'
$ws.Range('D41').ClearFormats()
$ws.Range('E41').Value = '  +1.29%  '
$ws.Range('D42').Value = '1.110.77'
$ws.Range('E42').Value = '  -3.87%  '
$ws.Range('B43').Value = 'FraxShare'
$ws.Range('C43').Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '6.094'
$ws.Range('D43').ClearFormats()
$ws.Range('E43').Value = '  +3.49%  '
$ws.Range('B44').Value = 'Aave'
$ws.Range('C44').Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '73.21'
$ws.Range('D44').ClearFormats()
$ws.Range('E44').Value = '  -0.81%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.8513'
$ws.Range('D45').ClearFormats()
$ws.Range('E45').Value = '  +0.32%  '
$ws.Range('E46').Value = '  -0.04%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '102.36'
$ws.Range('D47').ClearFormats()
$ws.Range('E47').Value = '  -0.22%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '1.870'
$ws.Range('D48').ClearFormats()
$ws.Range('E48').Value = '  -0.64%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '7.632'
$ws.Range('D49').ClearFormats()
$ws.Range('E49').Value = '  +2.52%  '
$ws.Range('D50').Value = '2.092.15'
$ws.Range('E50').Value = '  +2.48%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '2.986'
$ws.Range('D51').ClearFormats()
$ws.Range('E51').Value = '  -0.65%  '
